# Fix shark double counts in the Area Statistics / Sharks summary tables.
# Updates both the "Status by Landings (Area)" sheet and the
# "Status by Landings (Tier)" sheet, which both contain the same
# per-Sharks landings figures that were corrected for double counting.

$wb = $excel.ActiveWorkbook

# --- Sheet: "Status by Landings (Area)" ---
$wsArea = $wb.Worksheets.Item("Status by Landings (Area)")

$wsArea.Range("C3").Value  = 0.02118677
$wsArea.Range("C4").Value  = 0.00776831
$wsArea.Range("C5").Value  = 0.04865666
$wsArea.Range("C6").Value  = 0.00776831
$wsArea.Range("C7").Value  = 48.68392486517937
$wsArea.Range("C8").Value  = 37.54857113791996
$wsArea.Range("C9").Value  = 13.76750399690066
$wsArea.Range("C10").Value = 86.23249600309933
$wsArea.Range("C11").Value = 13.76750399690066

# --- Sheet: "Status by Landings (Tier)" ---
$wsTier = $wb.Worksheets.Item("Status by Landings (Tier)")

# Row 4: "Sharks"
$wsTier.Range("C4").Value = 0.02118677
$wsTier.Range("D4").Value = 0.00776831
$wsTier.Range("E4").Value = 0.04865666
$wsTier.Range("F4").Value = 0.00776831
$wsTier.Range("G4").Value = 48.68392486517937
$wsTier.Range("H4").Value = 37.54857113791996
$wsTier.Range("I4").Value = 13.76750399690066
$wsTier.Range("J4").Value = 86.23249600309933
$wsTier.Range("K4").Value = 13.76750399690066

# Row 5: "Global"
$wsTier.Range("C5").Value = 0.02118677
$wsTier.Range("D5").Value = 0.00776831
$wsTier.Range("E5").Value = 0.04865666
$wsTier.Range("F5").Value = 0.00776831
$wsTier.Range("G5").Value = 48.68392486517937
$wsTier.Range("H5").Value = 37.54857113791996
$wsTier.Range("I5").Value = 13.76750399690066
$wsTier.Range("J5").Value = 86.23249600309933
$wsTier.Range("K5").Value = 13.76750399690066
